$wb = $excel.ActiveWorkbook

# The handoff transform failed, so for both locale sheets (zh-cn, de-de)
# the first data row (row 2, the source .md file) needs to be updated to
# reflect the failure: status becomes "Handoff transform failed", the
# latest handoff file is cleared, the latest handoff datetime is reset to
# the empty/default datetime, and the handoff reason becomes "Ignored".

$sheetNames = @("zh-cn", "de-de")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # B2: Status -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # C2: Latest Handoff File -> cleared (also drops the hyperlink)
    $ws.Range("C2").ClearContents()

    # D2: Latest Handoff Datetime -> reset to default/empty datetime
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # H2: Handoff Reason -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}

# The Overview sheet shows the same "Status" value via the shared string
# table, so it needs to be brought in sync too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"
